# Add "2022-Q1" sheet (between "2021-Q4" and "总计") with one fund row,
# and prepend a "2022-Q1" summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook
$q4 = $wb.Worksheets.Item("2021-Q4")

# Remember the original "总计" rows before we touch anything.
$oldTotal = $wb.Worksheets.Item("总计")
$oldDates  = @()
$oldCounts = @()
$oldValues = @()
for ($r = 2; $r -le 5; $r++) {
    $oldDates  += $oldTotal.Cells.Item($r, 2).Value()
    $oldCounts += $oldTotal.Cells.Item($r, 3).Value()
    $oldValues += $oldTotal.Cells.Item($r, 4).Value()
}

# ---------------------------------------------------------------------
# 0) Delete the old "总计" sheet so its sheetId slot is freed, then
#    recreate the sheets in order so the new "2022-Q1" sheet reuses
#    that freed id and "总计" gets the next one (matches a from-scratch
#    export where sheetId just follows tab order: ...,4,5,6).
# ---------------------------------------------------------------------
$oldTotal.Delete() | Out-Null

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" worksheet right after "2021-Q4"
# ---------------------------------------------------------------------
$new = $wb.Worksheets.Add($null, $q4)
$new.Name = "2022-Q1"
$new.Outline.SummaryRow = 1
$new.Outline.SummaryColumn = 1
$new.PageSetup.LeftMargin = 54
$new.PageSetup.RightMargin = 54
$new.PageSetup.TopMargin = 72
$new.PageSetup.BottomMargin = 72
$new.PageSetup.HeaderMargin = 36
$new.PageSetup.FooterMargin = 36

# --- header row (B1:H1), styled like the other quarter sheets ---
$q4.Range("B1:H1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$new.Range("B1").Value() = "基金代码"
$new.Range("C1").Value() = "基金名称"
$new.Range("D1").Value() = "基金规模"
$new.Range("E1").Value() = "股票总仓位"
$new.Range("F1").Value() = "仓位占比"
$new.Range("G1").Value() = "持有市值(亿元)"
$new.Range("H1").Value() = "仓位排名"

# --- row-index cell (A2), styled like the other quarter sheets ---
$q4.Range("A2").Copy()
$new.Range("A2").PasteSpecial(-4122)
$new.Range("A2").Value() = 0

# --- data row 2: fund 000308 ---
# Columns B-G are stored as plain text (matches the other quarter
# sheets), H is numeric.
$new.Range("B2:G2").NumberFormat = "@"
$new.Range("B2").Value() = "000308"
$new.Range("C2").Value() = "建信创新中国混合"
$new.Range("D2").Value() = "3.11"
$new.Range("E2").Value() = "84.50"
$new.Range("F2").Value() = "3.02"
$new.Range("G2").Value() = "0.0939"
$new.Range("H2").Value() = 10

# drop the "@" quote-prefix flag so the cells end up with the default
# (unstyled) look, same as the sibling quarter sheets
$q4.Range("C2").Copy()
$new.Range("B2:G2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Recreate "总计" right after "2022-Q1" with the new summary row on
#    top, followed by the original rows (shifted down by one).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $new)
$total.Name = "总计"
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

$q4.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Range("B1").Value() = "日期"
$total.Range("C1").Value() = "持有数量(只)"
$total.Range("D1").Value() = "持有市值(亿元)"

$q4.Range("A2").Copy()
$total.Range("A2:A6").PasteSpecial(-4122)

$total.Range("A2").Value() = 0
$total.Range("B2").Value() = "2022-Q1"
$total.Range("C2").Value() = 1
$total.Range("D2").Value() = 0.09

for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 3
    $total.Cells.Item($r, 1).Value() = $i + 1
    $total.Cells.Item($r, 2).Value() = $oldDates[$i]
    $total.Cells.Item($r, 3).Value() = $oldCounts[$i]
    $total.Cells.Item($r, 4).Value() = $oldValues[$i]
}
